$wb = $excel.ActiveWorkbook

# "About" sheet gets a new state label and an updated date stamp
$ws = $wb.Worksheets.Item("About")

# New cell B1: "California" (state label next to the date)
$ws.Range("B1").Value = "California"

# Updated date in C1 (was 4/21/2021 -> now 11/9/2021)
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
